$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 926.2
$ws.Range("I4").Value = 658
$ws.Range("J4").Value = 1999
$ws.Range("K4").Value = 658
$ws.Range("L4").Value = 1999
$ws.Range("M4").Value = -544
$ws.Range("N4").Value = -2227
$ws.Range("H33").Value = 251.6
$ws.Range("I33").Value = 211
$ws.Range("K33").Value = 211
$ws.Range("M33").Value = 18
$ws.Range("H40").Value = 1699.7587
$ws.Range("I40").Value = 1711.2693
$ws.Range("K40").Value = 1711.2693
$ws.Range("M40").Value = -1536.2693
$ws.Range("H58").Value = 3102.6667
$ws.Range("I58").Value = 302.33334
$ws.Range("J58").Value = 4502.8335
$ws.Range("K58").Value = 907.0000200000001
$ws.Range("L58").Value = 13508.5005
$ws.Range("M58").Value = -757.0000200000001
$ws.Range("N58").Value = -13808.5005
$ws.Range("H88").Value = 1100
$ws.Range("J88").Value = 1100
$ws.Range("L88").Value = 1100
$ws.Range("N88").Value = -1912
$ws.Range("H91").Value = 1100
$ws.Range("J91").Value = 1100
$ws.Range("L91").Value = 1100
$ws.Range("N91").Value = -3908
$ws.Range("H112").Value = 3351.7222
$ws.Range("I112").Value = 2795
$ws.Range("J112").Value = 3384.4707
$ws.Range("K112").Value = 8385
$ws.Range("L112").Value = 10153.4121
$ws.Range("M112").Value = -7277
$ws.Range("N112").Value = -12369.4121
$ws.Range("H132").Value = 2861.5
$ws.Range("I132").Value = 2751.875
$ws.Range("K132").Value = 8255.625
$ws.Range("M132").Value = -5725.625
$ws.Range("H138").Value = 4157.027
$ws.Range("I138").Value = 3044.5
$ws.Range("J138").Value = 4691.04
$ws.Range("K138").Value = 9133.5
$ws.Range("L138").Value = 14073.12
$ws.Range("M138").Value = -3993.5
$ws.Range("N138").Value = -24353.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 9559
$ws.Range("I46").Value = 9750
$ws.Range("K46").Value = 9750
$ws.Range("M46").Value = -9431
$ws.Range("H102").Value = 659.8
$ws.Range("I102").Value = 599.75
$ws.Range("K102").Value = 599.75
$ws.Range("M102").Value = 1022.25
$ws.Range("H110").Value = 2623.8333
$ws.Range("I110").Value = 2623.8333
$ws.Range("K110").Value = 2623.8333
$ws.Range("M110").Value = -578.8332999999998
$ws.Range("H132").Value = 2195.75
$ws.Range("J132").Value = 2199
$ws.Range("L132").Value = 6597
$ws.Range("N132").Value = -11657

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5081.5
$ws.Range("I20").Value = 5267.1665
$ws.Range("K20").Value = 5267.1665
$ws.Range("M20").Value = -5020.1665
$ws.Range("H94").Value = 2257.2666
$ws.Range("I94").Value = 2061.5715
$ws.Range("K94").Value = 2061.5715
$ws.Range("M94").Value = -1610.5715
$ws.Range("H134").Value = 7515.25
$ws.Range("I134").Value = 7743.909
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 23231.727
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -20696.727
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 316.58334
$ws.Range("I7").Value = 89.7
$ws.Range("K7").Value = 89.7
$ws.Range("M7").Value = 23.3
$ws.Range("H62").Value = 2517.25
$ws.Range("I62").Value = 2485.5
$ws.Range("K62").Value = 2485.5
$ws.Range("M62").Value = -1861.5
$ws.Range("H65").Value = 2517.25
$ws.Range("I65").Value = 2485.5
$ws.Range("K65").Value = 12427.5
$ws.Range("M65").Value = -9307.5
$ws.Range("H99").Value = 4403.3335
$ws.Range("H126").Value = 4403.3335
$ws.Range("H132").Value = 8229.75
$ws.Range("I132").Value = 8229.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 24689.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -22159.25
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 2671.2856
$ws.Range("I134").Value = 2671.2856
$ws.Range("K134").Value = 8013.8568
$ws.Range("M134").Value = -5478.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 9233070
$ws.Range("K4").Value = 27699210
$ws.Range("M4").Value = -27699098
$ws.Range("H12").Value = 63.5
$ws.Range("J12").Value = 88
$ws.Range("L12").Value = 264
$ws.Range("N12").Value = -610
$ws.Range("H22").Value = 1379.6666
$ws.Range("I22").Value = 1282.6666
$ws.Range("J22").Value = 1476.6666
$ws.Range("K22").Value = 3847.9998
$ws.Range("L22").Value = 4429.9998
$ws.Range("M22").Value = -3678.9998
$ws.Range("N22").Value = -4767.9998
$ws.Range("H27").Value = 1379.6666
$ws.Range("I27").Value = 1282.6666
$ws.Range("J27").Value = 1476.6666
$ws.Range("K27").Value = 3847.9998
$ws.Range("L27").Value = 4429.9998
$ws.Range("M27").Value = -3745.9998
$ws.Range("N27").Value = -4633.9998
$ws.Range("H34").Value = 289.2
$ws.Range("J34").Value = 311.5
$ws.Range("L34").Value = 934.5
$ws.Range("N34").Value = -1102.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.611115
$ws.Range("J2").Value = 138.75
$ws.Range("L2").Value = 138.75
$ws.Range("N2").Value = -364.75
$ws.Range("H70").Value = 1496
$ws.Range("I70").Value = 1496
$ws.Range("K70").Value = 1496
$ws.Range("M70").Value = -1226
$ws.Range("H73").Value = 1496
$ws.Range("I73").Value = 1496
$ws.Range("K73").Value = 1496
$ws.Range("M73").Value = -560
$ws.Range("H80").Value = 5586.75
$ws.Range("I80").Value = 3733.3333
$ws.Range("J80").Value = 6698.8
$ws.Range("K80").Value = 3733.3333
$ws.Range("L80").Value = 6698.8
$ws.Range("M80").Value = -2735.3333
$ws.Range("N80").Value = -8694.799999999999
$ws.Range("H83").Value = 5586.75
$ws.Range("I83").Value = 3733.3333
$ws.Range("J83").Value = 6698.8
$ws.Range("K83").Value = 18666.6665
$ws.Range("L83").Value = 33494
$ws.Range("M83").Value = -13674.6665
$ws.Range("N83").Value = -43478
$ws.Range("H126").Value = 5022.769
$ws.Range("I126").Value = 3473.75
$ws.Range("J126").Value = 5711.222
$ws.Range("K126").Value = 10421.25
$ws.Range("L126").Value = 17133.666
$ws.Range("M126").Value = -7951.25
$ws.Range("N126").Value = -22073.666
$ws.Range("H132").Value = 3662.0833
$ws.Range("I132").Value = 3662.0833
$ws.Range("K132").Value = 10986.2499
$ws.Range("M132").Value = -8456.249899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2500
$ws.Range("I32").Value = 2500
$ws.Range("K32").Value = 2500
$ws.Range("M32").Value = -2183
$ws.Range("H46").Value = 1549.4445
$ws.Range("I46").Value = 1616.1666
$ws.Range("J46").Value = 1416
$ws.Range("K46").Value = 1616.1666
$ws.Range("L46").Value = 1416
$ws.Range("M46").Value = -1428.1666
$ws.Range("N46").Value = -1792
$ws.Range("H68").Value = 2930.3333
$ws.Range("J68").Value = 3048.5
$ws.Range("L68").Value = 3048.5
$ws.Range("N68").Value = -4546.5
$ws.Range("H71").Value = 2930.3333
$ws.Range("J71").Value = 3048.5
$ws.Range("L71").Value = 15242.5
$ws.Range("N71").Value = -22730.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14662.333
$ws.Range("I62").Value = 14662.333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 14662.333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -14038.333
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 14662.333
$ws.Range("I65").Value = 14662.333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 73311.66500000001
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -70191.66500000001
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 20244.666
$ws.Range("I74").Value = 19898
$ws.Range("K74").Value = 19898
$ws.Range("M74").Value = -18962
$ws.Range("H77").Value = 20244.666
$ws.Range("I77").Value = 19898
$ws.Range("K77").Value = 59694
$ws.Range("M77").Value = -55014
$ws.Range("H113").Value = 442.8889
$ws.Range("I113").Value = 414.5
$ws.Range("K113").Value = 1243.5
$ws.Range("M113").Value = 926.5
$ws.Range("H136").Value = 3868.5264
$ws.Range("I136").Value = 3914.4
$ws.Range("J136").Value = 3696.5
$ws.Range("K136").Value = 11743.2
$ws.Range("L136").Value = 11089.5
$ws.Range("M136").Value = -9193.200000000001
$ws.Range("N136").Value = -16189.5
